$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 34, shifting existing rows 34-131 down to 35-132
$ws.Rows.Item(34).Insert()

# Populate the newly inserted row 34 with the new data record
$ws.Cells.Item(34, 1).Value = 9
$ws.Cells.Item(34, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(34, 3).Value = "Metropolitana"
$ws.Cells.Item(34, 4).Value = 44414
$ws.Cells.Item(34, 5).Value = 13
$ws.Cells.Item(34, 6).Value = 300000001
$ws.Cells.Item(34, 7).Value = "Rabanito"
$ws.Cells.Item(34, 8).Value = "Sin especificar"
$ws.Cells.Item(34, 9).Value = "Primera"
$ws.Cells.Item(34, 10).Value = 7900
$ws.Cells.Item(34, 11).Value = 2500
$ws.Cells.Item(34, 12).Value = 3000
$ws.Cells.Item(34, 13).Value = 2750
$ws.Cells.Item(34, 14).Value = "`$/cien unidades (volumen en unidades)"
$ws.Cells.Item(34, 15).Value = "Provincia de Chacabuco"
$ws.Cells.Item(34, 16).Value = 28
$ws.Cells.Item(34, 17).Value = 100
$ws.Cells.Item(34, 18).Value = "Hortaliza"

# Copy the date style (numFmt) from row 35 (the shifted former row 34) into new row 34, column D
$ws.Cells.Item(35, 4).Copy()
$ws.Cells.Item(34, 4).PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(34, 4).Value = 44414
